$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 130870819
$ws.Range("B2").Value = 80348
$ws.Range("D2").Value = "NT"
$ws.Range("E2").Value = 6458
$ws.Range("F2").Value = "Lunglav"
$ws.Range("G2").Value = "Lobaria pulmonaria"
$ws.Range("H2").Value = "(L.) Hoffm."
$ws.Range("K2").Value = ""
$ws.Range("L2").Value = ""
$ws.Range("M2").Value = ""
$ws.Range("N2").Value = ""
$ws.Range("Q2").Value = 583057
$ws.Range("R2").Value = 6959583
$ws.Range("AC2").Value = ""

# Row 3
$ws.Range("A3").Value = 130870805
$ws.Range("B3").Value = 57884
$ws.Range("D3").Value = "NT"
$ws.Range("E3").Value = 100109
$ws.Range("F3").Value = "Tretåig hackspett"
$ws.Range("G3").Value = "Picoides tridactylus"
$ws.Range("H3").Value = "(Linnaeus, 1758)"
$ws.Range("K3").Value = ""
$ws.Range("L3").Value = ""
$ws.Range("M3").Value = "färska spår"
$ws.Range("N3").Value = ""
$ws.Range("Q3").Value = 582703
$ws.Range("R3").Value = 6959774
$ws.Range("AC3").Value = "Färska ringhack på tall"

# Row 4
$ws.Range("A4").Value = 130870796
$ws.Range("B4").Value = 57884
$ws.Range("D4").Value = "NT"
$ws.Range("E4").Value = 100109
$ws.Range("F4").Value = "Tretåig hackspett"
$ws.Range("G4").Value = "Picoides tridactylus"
$ws.Range("H4").Value = "(Linnaeus, 1758)"
$ws.Range("K4").Value = ""
$ws.Range("L4").Value = ""
$ws.Range("M4").Value = "äldre spår"
$ws.Range("N4").Value = ""
$ws.Range("Q4").Value = 582707
$ws.Range("R4").Value = 6959492
$ws.Range("AC4").Value = "Äldre ringhack på tall"

# Row 10
$ws.Range("A10").Value = 130870826
$ws.Range("B10").Value = 79243
$ws.Range("D10").Value = "NT"
$ws.Range("E10").Value = 6425
$ws.Range("F10").Value = "Garnlav"
$ws.Range("G10").Value = "Alectoria sarmentosa"
$ws.Range("H10").Value = "(Ach.) Ach."
$ws.Range("K10").Value = ""
$ws.Range("L10").Value = ""
$ws.Range("M10").Value = ""
$ws.Range("N10").Value = ""
$ws.Range("Q10").Value = 582992
$ws.Range("R10").Value = 6959624
$ws.Range("AC10").Value = ""

# Row 11
$ws.Range("A11").Value = 130870828
$ws.Range("B11").Value = 79243
$ws.Range("D11").Value = "NT"
$ws.Range("E11").Value = 6425
$ws.Range("F11").Value = "Garnlav"
$ws.Range("G11").Value = "Alectoria sarmentosa"
$ws.Range("H11").Value = "(Ach.) Ach."
$ws.Range("K11").Value = ""
$ws.Range("L11").Value = ""
$ws.Range("M11").Value = ""
$ws.Range("N11").Value = ""
$ws.Range("Q11").Value = 583216
$ws.Range("R11").Value = 6959386
$ws.Range("AC11").Value = ""

# Row 15
$ws.Range("A15").Value = 130870831
$ws.Range("B15").Value = 83089
$ws.Range("D15").Value = "NT"
$ws.Range("E15").Value = 1312
$ws.Range("F15").Value = "Gammelgransskål"
$ws.Range("G15").Value = "Pseudographis pinicola"
$ws.Range("H15").Value = "(Nyl.) Rehm"
$ws.Range("K15").Value = ""
$ws.Range("L15").Value = ""
$ws.Range("M15").Value = ""
$ws.Range("N15").Value = ""
$ws.Range("Q15").Value = 583209
$ws.Range("R15").Value = 6959416
$ws.Range("AC15").Value = ""

# Row 16
$ws.Range("A16").Value = 130870818
$ws.Range("B16").Value = 92267
$ws.Range("D16").Value = "VU"
$ws.Range("E16").Value = 1209
$ws.Range("F16").Value = "Rynkskinn"
$ws.Range("G16").Value = "Hermanssonia centrifuga"
$ws.Range("H16").Value = "(P. Karst.) Zmitr."
$ws.Range("K16").Value = ""
$ws.Range("L16").Value = ""
$ws.Range("M16").Value = ""
$ws.Range("N16").Value = ""
$ws.Range("Q16").Value = 583241
$ws.Range("R16").Value = 6959405
$ws.Range("AC16").Value = ""

# Row 17
$ws.Range("A17").Value = 130870792
$ws.Range("B17").Value = 91808
$ws.Range("D17").Value = "NT"
$ws.Range("E17").Value = 1202
$ws.Range("F17").Value = "Ullticka"
$ws.Range("G17").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H17").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("K17").Value = ""
$ws.Range("L17").Value = ""
$ws.Range("M17").Value = ""
$ws.Range("N17").Value = ""
$ws.Range("Q17").Value = 583131
$ws.Range("R17").Value = 6959482
$ws.Range("AC17").Value = ""

# Row 18
$ws.Range("A18").Value = 130870823
$ws.Range("B18").Value = 79243
$ws.Range("D18").Value = "NT"
$ws.Range("E18").Value = 6425
$ws.Range("F18").Value = "Garnlav"
$ws.Range("G18").Value = "Alectoria sarmentosa"
$ws.Range("H18").Value = "(Ach.) Ach."
$ws.Range("K18").Value = ""
$ws.Range("L18").Value = ""
$ws.Range("M18").Value = ""
$ws.Range("N18").Value = ""
$ws.Range("Q18").Value = 582529
$ws.Range("R18").Value = 6959663
$ws.Range("AC18").Value = "Med apothecier"

# Row 19
$ws.Range("A19").Value = 130870816
$ws.Range("B19").Value = 80377
$ws.Range("D19").Value = "LC"
$ws.Range("E19").Value = 6462
$ws.Range("F19").Value = "Stuplav"
$ws.Range("G19").Value = "Nephroma bellum"
$ws.Range("H19").Value = "(Spreng.) Tuck."
$ws.Range("K19").Value = ""
$ws.Range("L19").Value = ""
$ws.Range("M19").Value = ""
$ws.Range("N19").Value = ""
$ws.Range("Q19").Value = 582711
$ws.Range("R19").Value = 6959664
$ws.Range("AC19").Value = ""

# Row 20
$ws.Range("A20").Value = 130870795
$ws.Range("B20").Value = 57884
$ws.Range("D20").Value = "NT"
$ws.Range("E20").Value = 100109
$ws.Range("F20").Value = "Tretåig hackspett"
$ws.Range("G20").Value = "Picoides tridactylus"
$ws.Range("H20").Value = "(Linnaeus, 1758)"
$ws.Range("K20").Value = ""
$ws.Range("L20").Value = ""
$ws.Range("M20").Value = "äldre spår"
$ws.Range("N20").Value = ""
$ws.Range("Q20").Value = 583098
$ws.Range("R20").Value = 6959481
$ws.Range("AC20").Value = "Äldre ringhack på tall"

# Row 21
$ws.Range("A21").Value = 130870798
$ws.Range("B21").Value = 57884
$ws.Range("D21").Value = "NT"
$ws.Range("E21").Value = 100109
$ws.Range("F21").Value = "Tretåig hackspett"
$ws.Range("G21").Value = "Picoides tridactylus"
$ws.Range("H21").Value = "(Linnaeus, 1758)"
$ws.Range("K21").Value = ""
$ws.Range("L21").Value = ""
$ws.Range("M21").Value = "färska spår"
$ws.Range("N21").Value = ""
$ws.Range("Q21").Value = 582557
$ws.Range("R21").Value = 6959519
$ws.Range("AC21").Value = "Färska ringhack på tall"

# Row 22
$ws.Range("A22").Value = 130870797
$ws.Range("B22").Value = 57884
$ws.Range("D22").Value = "NT"
$ws.Range("E22").Value = 100109
$ws.Range("F22").Value = "Tretåig hackspett"
$ws.Range("G22").Value = "Picoides tridactylus"
$ws.Range("H22").Value = "(Linnaeus, 1758)"
$ws.Range("K22").Value = ""
$ws.Range("L22").Value = ""
$ws.Range("M22").Value = "färska spår"
$ws.Range("N22").Value = ""
$ws.Range("Q22").Value = 582602
$ws.Range("R22").Value = 6959501
$ws.Range("AC22").Value = "Färska ringhack på tall"

# Row 23
$ws.Range("A23").Value = 130870817
$ws.Range("B23").Value = 91819
$ws.Range("D23").Value = "LC"
$ws.Range("E23").Value = 1205
$ws.Range("F23").Value = "Stor aspticka"
$ws.Range("G23").Value = "Phellinus populicola"
$ws.Range("H23").Value = "Niemelä"
$ws.Range("K23").Value = ""
$ws.Range("L23").Value = ""
$ws.Range("M23").Value = ""
$ws.Range("N23").Value = ""
$ws.Range("Q23").Value = 582663
$ws.Range("R23").Value = 6959537
$ws.Range("AC23").Value = ""

# Row 27
$ws.Range("A27").Value = 130870825
$ws.Range("B27").Value = 79243
$ws.Range("D27").Value = "NT"
$ws.Range("E27").Value = 6425
$ws.Range("F27").Value = "Garnlav"
$ws.Range("G27").Value = "Alectoria sarmentosa"
$ws.Range("H27").Value = "(Ach.) Ach."
$ws.Range("K27").Value = ""
$ws.Range("L27").Value = ""
$ws.Range("M27").Value = ""
$ws.Range("N27").Value = ""
$ws.Range("Q27").Value = 582799
$ws.Range("R27").Value = 6959667
$ws.Range("AC27").Value = ""

# Row 28
$ws.Range("A28").Value = 130870824
$ws.Range("B28").Value = 79243
$ws.Range("D28").Value = "NT"
$ws.Range("E28").Value = 6425
$ws.Range("F28").Value = "Garnlav"
$ws.Range("G28").Value = "Alectoria sarmentosa"
$ws.Range("H28").Value = "(Ach.) Ach."
$ws.Range("K28").Value = ""
$ws.Range("L28").Value = ""
$ws.Range("M28").Value = ""
$ws.Range("N28").Value = ""
$ws.Range("Q28").Value = 582685
$ws.Range("R28").Value = 6959786
$ws.Range("AC28").Value = ""

# Row 29
$ws.Range("A29").Value = 130870815
$ws.Range("B29").Value = 57884
$ws.Range("D29").Value = "NT"
$ws.Range("E29").Value = 100109
$ws.Range("F29").Value = "Tretåig hackspett"
$ws.Range("G29").Value = "Picoides tridactylus"
$ws.Range("H29").Value = "(Linnaeus, 1758)"
$ws.Range("K29").Value = ""
$ws.Range("L29").Value = ""
$ws.Range("M29").Value = "färska spår"
$ws.Range("N29").Value = ""
$ws.Range("Q29").Value = 583170
$ws.Range("R29").Value = 6959447
$ws.Range("AC29").Value = "Färska och äldre ringhack på tall"

# Row 30
$ws.Range("A30").Value = 130870827
$ws.Range("B30").Value = 79243
$ws.Range("D30").Value = "NT"
$ws.Range("E30").Value = 6425
$ws.Range("F30").Value = "Garnlav"
$ws.Range("G30").Value = "Alectoria sarmentosa"
$ws.Range("H30").Value = "(Ach.) Ach."
$ws.Range("K30").Value = ""
$ws.Range("L30").Value = ""
$ws.Range("M30").Value = ""
$ws.Range("N30").Value = ""
$ws.Range("Q30").Value = 583142
$ws.Range("R30").Value = 6959494
$ws.Range("AC30").Value = ""

# Row 33
$ws.Range("A33").Value = 130870794
$ws.Range("B33").Value = 91808
$ws.Range("D33").Value = "NT"
$ws.Range("E33").Value = 1202
$ws.Range("F33").Value = "Ullticka"
$ws.Range("G33").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H33").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("K33").Value = ""
$ws.Range("L33").Value = ""
$ws.Range("M33").Value = ""
$ws.Range("N33").Value = ""
$ws.Range("Q33").Value = 583237
$ws.Range("R33").Value = 6959408
$ws.Range("AC33").Value = ""

# Row 34
$ws.Range("A34").Value = 130870832
$ws.Range("B34").Value = 80383
$ws.Range("D34").Value = "LC"
$ws.Range("E34").Value = 6463
$ws.Range("F34").Value = "Bårdlav"
$ws.Range("G34").Value = "Nephroma parile"
$ws.Range("H34").Value = "(Ach.) Ach."
$ws.Range("K34").Value = ""
$ws.Range("L34").Value = ""
$ws.Range("M34").Value = ""
$ws.Range("N34").Value = ""
$ws.Range("Q34").Value = 583054
$ws.Range("R34").Value = 6959568
$ws.Range("AC34").Value = ""

# Row 35
$ws.Range("A35").Value = 130870811
$ws.Range("B35").Value = 57884
$ws.Range("D35").Value = "NT"
$ws.Range("E35").Value = 100109
$ws.Range("F35").Value = "Tretåig hackspett"
$ws.Range("G35").Value = "Picoides tridactylus"
$ws.Range("H35").Value = "(Linnaeus, 1758)"
$ws.Range("K35").Value = ""
$ws.Range("L35").Value = ""
$ws.Range("M35").Value = "äldre spår"
$ws.Range("N35").Value = ""
$ws.Range("Q35").Value = 582879
$ws.Range("R35").Value = 6959670
$ws.Range("AC35").Value = "Äldre ringhack på tall"

# Row 36
$ws.Range("A36").Value = 130870812
$ws.Range("B36").Value = 57884
$ws.Range("D36").Value = "NT"
$ws.Range("E36").Value = 100109
$ws.Range("F36").Value = "Tretåig hackspett"
$ws.Range("G36").Value = "Picoides tridactylus"
$ws.Range("H36").Value = "(Linnaeus, 1758)"
$ws.Range("K36").Value = ""
$ws.Range("L36").Value = ""
$ws.Range("M36").Value = "färska spår"
$ws.Range("N36").Value = ""
$ws.Range("Q36").Value = 582898
$ws.Range("R36").Value = 6959678
$ws.Range("AC36").Value = "Färska och äldre ringhack på tall"

# Row 37
$ws.Range("A37").Value = 130870810
$ws.Range("B37").Value = 57884
$ws.Range("D37").Value = "NT"
$ws.Range("E37").Value = 100109
$ws.Range("F37").Value = "Tretåig hackspett"
$ws.Range("G37").Value = "Picoides tridactylus"
$ws.Range("H37").Value = "(Linnaeus, 1758)"
$ws.Range("K37").Value = ""
$ws.Range("L37").Value = ""
$ws.Range("M37").Value = "äldre spår"
$ws.Range("N37").Value = ""
$ws.Range("Q37").Value = 582825
$ws.Range("R37").Value = 6959676
$ws.Range("AC37").Value = "Äldre ringhack på tall"

# Row 38
$ws.Range("A38").Value = 130870808
$ws.Range("B38").Value = 57884
$ws.Range("D38").Value = "NT"
$ws.Range("E38").Value = 100109
$ws.Range("F38").Value = "Tretåig hackspett"
$ws.Range("G38").Value = "Picoides tridactylus"
$ws.Range("H38").Value = "(Linnaeus, 1758)"
$ws.Range("K38").Value = ""
$ws.Range("L38").Value = ""
$ws.Range("M38").Value = "färska spår"
$ws.Range("N38").Value = ""
$ws.Range("Q38").Value = 582781
$ws.Range("R38").Value = 6959717
$ws.Range("AC38").Value = "Färska ringhack på tall"

# Row 39
$ws.Range("A39").Value = 130870791
$ws.Range("B39").Value = 91808
$ws.Range("D39").Value = "NT"
$ws.Range("E39").Value = 1202
$ws.Range("F39").Value = "Ullticka"
$ws.Range("G39").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H39").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("K39").Value = ""
$ws.Range("L39").Value = ""
$ws.Range("M39").Value = ""
$ws.Range("N39").Value = ""
$ws.Range("Q39").Value = 582769
$ws.Range("R39").Value = 6959516
$ws.Range("AC39").Value = ""
